$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the used range extent
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Replace "12.5.18" with "12.05.18" in column D (libraryDate), for all data rows
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($val -eq "12.5.18") {
        $cell.Value = "12.05.18"
    }
}

# Update the active selection to match the saved view state (E27)
$ws.Range("E27").Select()
